$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "98-44="
$tbl.Cell(1, 2).Range.Text = "33-15="
$tbl.Cell(1, 3).Range.Text = "20+42="
$tbl.Cell(1, 4).Range.Text = "14+30="
$tbl.Cell(1, 5).Range.Text = "41+9="

$tbl.Cell(2, 1).Range.Text = "19+1="
$tbl.Cell(2, 2).Range.Text = "37+45="
$tbl.Cell(2, 3).Range.Text = "48+39="
$tbl.Cell(2, 4).Range.Text = "97-69="
$tbl.Cell(2, 5).Range.Text = "17+69="

$tbl.Cell(3, 1).Range.Text = "91-81="
$tbl.Cell(3, 2).Range.Text = "84-53="
$tbl.Cell(3, 3).Range.Text = "80-47="
$tbl.Cell(3, 4).Range.Text = "36+50="
$tbl.Cell(3, 5).Range.Text = "51+18="

$tbl.Cell(4, 1).Range.Text = "44-6="
$tbl.Cell(4, 2).Range.Text = "43-23="
$tbl.Cell(4, 3).Range.Text = "21+5="
$tbl.Cell(4, 4).Range.Text = "35+14="
$tbl.Cell(4, 5).Range.Text = "64-63="

$tbl.Cell(5, 1).Range.Text = "49+37="
$tbl.Cell(5, 2).Range.Text = "38+45="
$tbl.Cell(5, 3).Range.Text = "65+5="
$tbl.Cell(5, 4).Range.Text = "57-24="
$tbl.Cell(5, 5).Range.Text = "37-14="

$tbl.Cell(6, 1).Range.Text = "30+27="
$tbl.Cell(6, 2).Range.Text = "18+22="
$tbl.Cell(6, 3).Range.Text = "62-48="
$tbl.Cell(6, 4).Range.Text = "54+43="
$tbl.Cell(6, 5).Range.Text = "14+18="

$tbl.Cell(7, 1).Range.Text = "91-75="
$tbl.Cell(7, 2).Range.Text = "10+56="
$tbl.Cell(7, 3).Range.Text = "35-29="
$tbl.Cell(7, 4).Range.Text = "61-53="
$tbl.Cell(7, 5).Range.Text = "15+41="

$tbl.Cell(8, 1).Range.Text = "15-10="
$tbl.Cell(8, 2).Range.Text = "75-21="
$tbl.Cell(8, 3).Range.Text = "19+73="
$tbl.Cell(8, 4).Range.Text = "77-16="
$tbl.Cell(8, 5).Range.Text = "3+78="

$tbl.Cell(9, 1).Range.Text = "88-2="
$tbl.Cell(9, 2).Range.Text = "81-32="
$tbl.Cell(9, 3).Range.Text = "1+40="
$tbl.Cell(9, 4).Range.Text = "43+18="
$tbl.Cell(9, 5).Range.Text = "10+34="

$tbl.Cell(10, 1).Range.Text = "56-45="
$tbl.Cell(10, 2).Range.Text = "5+63="
$tbl.Cell(10, 3).Range.Text = "91-8="
$tbl.Cell(10, 4).Range.Text = "20+60="
$tbl.Cell(10, 5).Range.Text = "86-10="

$tbl.Cell(11, 1).Range.Text = "35+63="
$tbl.Cell(11, 2).Range.Text = "66-38="
$tbl.Cell(11, 3).Range.Text = "3+38="
$tbl.Cell(11, 4).Range.Text = "88-19="
$tbl.Cell(11, 5).Range.Text = "34-27="

$tbl.Cell(12, 1).Range.Text = "61-26="
$tbl.Cell(12, 2).Range.Text = "2+75="
$tbl.Cell(12, 3).Range.Text = "18+31="
$tbl.Cell(12, 4).Range.Text = "13-7="
$tbl.Cell(12, 5).Range.Text = "68-8="

$tbl.Cell(13, 1).Range.Text = "55+32="
$tbl.Cell(13, 2).Range.Text = "28+47="
$tbl.Cell(13, 3).Range.Text = "59-50="
$tbl.Cell(13, 4).Range.Text = "54-15="
$tbl.Cell(13, 5).Range.Text = "8+91="

$tbl.Cell(14, 1).Range.Text = "68-38="
$tbl.Cell(14, 2).Range.Text = "49+23="
$tbl.Cell(14, 3).Range.Text = "52-45="
$tbl.Cell(14, 4).Range.Text = "97-69="
$tbl.Cell(14, 5).Range.Text = "41+24="

$tbl.Cell(15, 1).Range.Text = "84-56="
$tbl.Cell(15, 2).Range.Text = "59-7="
$tbl.Cell(15, 3).Range.Text = "88-29="
$tbl.Cell(15, 4).Range.Text = "63+4="
$tbl.Cell(15, 5).Range.Text = "11+59="

$tbl.Cell(16, 1).Range.Text = "63-62="
$tbl.Cell(16, 2).Range.Text = "42+15="
$tbl.Cell(16, 3).Range.Text = "34+58="
$tbl.Cell(16, 4).Range.Text = "6+62="
$tbl.Cell(16, 5).Range.Text = "16+48="

$tbl.Cell(17, 1).Range.Text = "58-17="
$tbl.Cell(17, 2).Range.Text = "90-42="
$tbl.Cell(17, 3).Range.Text = "80-8="
$tbl.Cell(17, 4).Range.Text = "55+31="
$tbl.Cell(17, 5).Range.Text = "61+13="

$tbl.Cell(18, 1).Range.Text = "54-12="
$tbl.Cell(18, 2).Range.Text = "72+16="
$tbl.Cell(18, 3).Range.Text = "4+36="
$tbl.Cell(18, 4).Range.Text = "31-3="
$tbl.Cell(18, 5).Range.Text = "25-3="

$tbl.Cell(19, 1).Range.Text = "55-6="
$tbl.Cell(19, 2).Range.Text = "54-22="
$tbl.Cell(19, 3).Range.Text = "97-32="
$tbl.Cell(19, 4).Range.Text = "5+28="
$tbl.Cell(19, 5).Range.Text = "39+56="

$tbl.Cell(20, 1).Range.Text = "95-9="
$tbl.Cell(20, 2).Range.Text = "72-56="
$tbl.Cell(20, 3).Range.Text = "17+45="
$tbl.Cell(20, 4).Range.Text = "23+58="
$tbl.Cell(20, 5).Range.Text = "34-19="

